$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "69-36=33","98-67=31","66+28=94","46+3=49","35+9=44",
    "39-8=31","15+66=81","10+9=19","36+60=96","4+22=26",
    "74-46=28","9+21=30","59+27=86","98-41=57","38+11=49",
    "23+15=38","38+23=61","26+52=78","62+24=86","24-17=7",
    "81-46=35","73-34=39","0+89=89","97-48=49","89-82=7",
    "89-78=11","3+57=60","86-31=55","84-56=28","71+27=98",
    "5+73=78","10+2=12","65-56=9","62+29=91","22-11=11",
    "31-17=14","65-4=61","81-68=13","99-75=24","56-53=3",
    "31+27=58","79+10=89","24+19=43","32-14=18","9+36=45",
    "17+65=82","28+43=71","25-14=11","99-62=37","93-39=54",
    "94-55=39","52+20=72","71-43=28","78+11=89","97-71=26",
    "13+76=89","84-43=41","51-17=34","22-0=22","56+43=99",
    "68-43=25","17+24=41","33-31=2","66-58=8","0+31=31",
    "66-20=46","54-33=21","31+37=68","69-54=15","24+2=26",
    "62-30=32","9+0=9","45+21=66","66-13=53","60-14=46",
    "13+11=24","9+67=76","29+33=62","15+18=33","97-83=14",
    "76-40=36","76+23=99","81-20=61","19+25=44","40+3=43",
    "40-22=18","98-48=50","30+28=58","88-12=76","7+77=84",
    "85-43=42","98-89=9","65+25=90","51-5=46","73+20=93",
    "68-62=6","13+59=72","41+47=88","89-24=65","96+1=97"
)

$rows = 20
$cols = 5
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $idx = ($r - 1) * $cols + ($c - 1)
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
    }
}

Write-Host "Done updating" $rows "rows x" $cols "cols"
